# Fix contact information missing from short resumes:
# Insert a new centered paragraph with contact info right after the
# name heading ("Dheeraj Chand").
$d = $word.ActiveDocument

$null = $d.Content.Find.Execute(
    "Dheeraj Chand",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "Dheeraj Chand^p202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX",
    2
)
